$wb = $excel.ActiveWorkbook

# "Level2_EconTest" sheet: B2 and C2 were 0, should become 1 (need >=3 non-zero
# criteria, fixing "Cannot work on less than 3 criteria" failure).
$wsTest = $wb.Worksheets.Item("Level2_EconTest")
$wsTest.Range("B2").Value = 1
$wsTest.Range("C2").Value = 1

# Make "Level2_EconTest" the active sheet/tab with C2 selected.
# (Previously "Level2_Econ" was the active tab; activating this sheet moves the
# tabSelected flag here and clears it from "Level2_Econ".)
$wsTest.Activate()
$wsTest.Range("C2").Select()
